$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company name for row 3
$ws.Range("B3").Value = "Ice Group ASA (OB:ICEGR)"

# Row 2 updates
$ws.Range("G2").Value = -0.14233907524932
$ws.Range("H2").Value = -0.14233907524932
$ws.Range("I2").Value = -0.242973708068903
$ws.Range("J2").Value = -0.242973708068903
$ws.Range("K2").Value = -113.6
$ws.Range("L2").Value = -0.514959202175884
$ws.Range("U2").Value = 119.9
$ws.Range("V2").Value = 0.3305762338020403
$ws.Range("W2").Value = 1.627507163323782
$ws.Range("X2").Value = 0.08765733785408866
$ws.Range("Y2").Value = 1.539849825469694
$ws.Range("Z2").Value = 0.5668036998972251
$ws.Range("AA2").Value = -0.1377183967112025
$ws.Range("AB2").Value = 0.04195107328433252
$ws.Range("AC2").Value = -0.179669469995535
$ws.Range("AD2").Value = 759.5
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 759.5
$ws.Range("AG2").Value = 639.6
$ws.Range("AH2").Value = 0.6767955801104972
$ws.Range("AI2").Value = 1.303862660944206
$ws.Range("AJ2").Value = 0.6381322957198444
$ws.Range("AK2").Value = 1.382619974059663
$ws.Range("AL2").Value = 56.9
$ws.Range("AM2").Value = 56.58
$ws.Range("AN2").Value = -59.3359375
$ws.Range("AO2").Value = -0.9420035149384887
$ws.Range("AP2").Value = -49.96875
$ws.Range("AQ2").Value = -0.9473312124425592

# Row 3 updates
$ws.Range("G3").Value = -0.14233907524932
$ws.Range("H3").Value = -0.14233907524932
$ws.Range("I3").Value = -0.242973708068903
$ws.Range("J3").Value = -0.242973708068903
$ws.Range("K3").Value = -113.6
$ws.Range("L3").Value = -0.514959202175884
$ws.Range("U3").Value = 119.9
$ws.Range("V3").Value = 0.3305762338020403
$ws.Range("W3").Value = 1.627507163323782
$ws.Range("X3").Value = 0.08765733785408866
$ws.Range("Y3").Value = 1.539849825469694
$ws.Range("Z3").Value = 0.5668036998972251
$ws.Range("AA3").Value = -0.1377183967112025
$ws.Range("AB3").Value = 0.04195107328433252
$ws.Range("AC3").Value = -0.179669469995535
$ws.Range("AD3").Value = 759.5
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 759.5
$ws.Range("AG3").Value = 639.6
$ws.Range("AH3").Value = 0.6767955801104972
$ws.Range("AI3").Value = 1.303862660944206
$ws.Range("AJ3").Value = 0.6381322957198444
$ws.Range("AK3").Value = 1.382619974059663
$ws.Range("AL3").Value = 56.9
$ws.Range("AM3").Value = 56.58
$ws.Range("AN3").Value = -59.3359375
$ws.Range("AO3").Value = -0.9420035149384887
$ws.Range("AP3").Value = -49.96875
$ws.Range("AQ3").Value = -0.9473312124425592
